$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) cells are stored as Text so numeric-looking strings
# (e.g. "1.00", "574.46") are not coerced into numbers - matches the
# original workbook, where column D is authored as inline/shared text.

# Row 2
$ws.Cells.Item(2, 4).Value = "65.525.63"
$ws.Cells.Item(2, 5).Value = "  -0.82%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.276.45"
$ws.Cells.Item(3, 5).Value = "  -1.15%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "574.46"
$ws.Cells.Item(5, 5).Value = "  +2.63%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "181.16"
$ws.Cells.Item(6, 5).Value = "  -3.51%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.23%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.273.27"
$ws.Cells.Item(8, 5).Value = "  -0.89%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.568"
$ws.Cells.Item(9, 5).Value = "  -2.86%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.174"
$ws.Cells.Item(10, 5).Value = "  -6.23%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.567"
$ws.Cells.Item(11, 5).Value = "  -3.01%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "45.91"
$ws.Cells.Item(12, 5).Value = "  -3.74%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000261"
$ws.Cells.Item(13, 5).Value = "  -3.56%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.798.30"
$ws.Cells.Item(14, 5).Value = "  -1.20%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "8.34"
$ws.Cells.Item(15, 5).Value = "  -3.42%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "607.75"
$ws.Cells.Item(16, 5).Value = "  -3.46%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "65.626.96"
$ws.Cells.Item(17, 5).Value = "  -0.68%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +0.22%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(19, 4).Value = "3.282.26"
$ws.Cells.Item(19, 5).Value = "  -0.28%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Chainlink"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "17.60"
$ws.Cells.Item(20, 5).Value = "  -3.11%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.82"
$ws.Cells.Item(21, 5).Value = "  -2.52%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.882"
$ws.Cells.Item(22, 5).Value = "  -3.02%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "18.14"
$ws.Cells.Item(23, 5).Value = "  -0.10%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.38%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "97.69"
$ws.Cells.Item(25, 5).Value = "  -4.73%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.93"
$ws.Cells.Item(26, 5).Value = "  -0.19%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.72%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.37"
$ws.Cells.Item(28, 5).Value = "  -1.98%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "30.57"
$ws.Cells.Item(29, 5).Value = "  +0.92%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.32"
$ws.Cells.Item(30, 5).Value = "  -4.16%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.42"
$ws.Cells.Item(31, 5).Value = "  +1.06%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.70"
$ws.Cells.Item(32, 5).Value = "  -8.87%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Bittensor"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "544.11"
$ws.Cells.Item(33, 5).Value = "  -1.41%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Cosmos"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "10.77"
$ws.Cells.Item(34, 5).Value = "  -2.81%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "3.792.93"
$ws.Cells.Item(35, 5).Value = "  -1.60%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.102"
$ws.Cells.Item(36, 5).Value = "  -3.05%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.998"
$ws.Cells.Item(37, 5).Value = "  -0.12%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "55.92"
$ws.Cells.Item(38, 5).Value = "  -2.89%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.86%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "ApeXProtocol"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.40"
$ws.Cells.Item(40, 5).Value = "  +4.92%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "32.29"
$ws.Cells.Item(41, 5).Value = "  -4.32%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.12"
$ws.Cells.Item(42, 5).Value = "  -5.40%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "0.0₃0675"
$ws.Cells.Item(43, 5).Value = "  -7.95%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.56"
$ws.Cells.Item(44, 5).Value = "  -4.54%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.327"
$ws.Cells.Item(45, 5).Value = "  -1.96%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0402"
$ws.Cells.Item(46, 5).Value = "  -4.33%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.02"
$ws.Cells.Item(47, 5).Value = "  -7.08%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.40%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.125"
$ws.Cells.Item(49, 5).Value = "  -2.81%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -4.66%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "127.68"
$ws.Cells.Item(51, 5).Value = "  +4.33%  "
